$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add the new "Same Tree" entry as row 21 on Sheet1
$ws1.Range("B21").Value = "Same Tree"
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 2
$ws1.Range("E21").Value = 38
$ws1.Range("F21").Value = 0.43
$ws1.Range("G21").Value = 16.5
$ws1.Range("H21").Value = 0.89
$ws1.Range("I21").Value = "https://leetcode.com/problems/same-tree/submissions/1061449405/"

# Update selections / active sheet to match the authored state
$ws2.Select()
$ws2.Range("E11").Select()
$ws1.Select()
$ws1.Range("B25").Select()
